# Apply updated cryptocurrency data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    $c = $ws.Range($range)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell "D2" "43.286.11"
Set-TextCell "E2" "  -4.75%  "
Set-TextCell "D3" "2.244.43"
Set-TextCell "E3" "  -5.69%  "
Set-TextCell "E4" "  -0.13%  "
Set-TextCell "D5" "321.35"
Set-TextCell "E5" "  +0.50%  "
Set-TextCell "D6" "101.56"
Set-TextCell "E6" "  -7.14%  "
Set-TextCell "D7" "0.586"
Set-TextCell "E7" "  -7.90%  "
Set-TextCell "E8" "  -0.18%  "
Set-TextCell "D9" "0.566"
Set-TextCell "E9" "  -8.30%  "
Set-TextCell "D10" "37.35"
Set-TextCell "E10" "  -9.29%  "
Set-TextCell "D11" "54.64"
Set-TextCell "E11" "  -2.78%  "
Set-TextCell "D12" "0.0831"
Set-TextCell "E12" "  -9.71%  "
Set-TextCell "D13" "7.76"
Set-TextCell "E13" "  -9.37%  "
Set-TextCell "E14" "  -0.53%  "
Set-TextCell "D15" "0.871"
Set-TextCell "E15" "  -11.88%  "
Set-TextCell "D16" "2.585.44"
Set-TextCell "E16" "  -5.70%  "
Set-TextCell "D17" "14.52"
Set-TextCell "E17" "  -6.48%  "
Set-TextCell "D18" "2.234.27"
Set-TextCell "E18" "  -6.09%  "
Set-TextCell "D19" "43.225.09"
Set-TextCell "E19" "  -4.76%  "
Set-TextCell "D20" "14.56"
Set-TextCell "E20" "  -7.08%  "
Set-TextCell "D21" "0.0₃0971"
Set-TextCell "E21" "  -8.72%  "
Set-TextCell "D22" "6.58"
Set-TextCell "E22" "  -10.44%  "
Set-TextCell "D23" "65.70"
Set-TextCell "E23" "  -10.45%  "
Set-TextCell "E24" "  -14.68%  "
Set-TextCell "D25" "239.13"
Set-TextCell "E25" "  -8.58%  "
Set-TextCell "D26" "2.18"
Set-TextCell "E26" "  -7.31%  "
Set-TextCell "E27" "  -0.05%  "
Set-TextCell "E28" "  +1.80%  "
Set-TextCell "B29" "Cosmos"
Set-TextCell "C29" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextCell "D29" "10.09"
Set-TextCell "E29" "  -10.70%  "
Set-TextCell "B30" "Toncoin"
Set-TextCell "C30" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextCell "D30" "2.25"
Set-TextCell "E30" "  +1.38%  "
Set-TextCell "D31" "6.42"
Set-TextCell "E31" "  -15.41%  "
Set-TextCell "D32" "36.30"
Set-TextCell "E32" "  -3.45%  "
Set-TextCell "D33" "0.0886"
Set-TextCell "E33" "  -7.17%  "
Set-TextCell "D34" "20.52"
Set-TextCell "E34" "  -8.67%  "
Set-TextCell "D35" "153.45"
Set-TextCell "E35" "  -8.59%  "
Set-TextCell "D36" "2.71"
Set-TextCell "E36" "  -6.58%  "
Set-TextCell "D37" "3.23"
Set-TextCell "E37" "  +8.19%  "
Set-TextCell "D38" "1.97"
Set-TextCell "E38" "  +0.04%  "
Set-TextCell "E39" "  -7.47%  "
Set-TextCell "D40" "4.47"
Set-TextCell "E40" "  -5.96%  "
Set-TextCell "E41" "  -10.75%  "
Set-TextCell "D42" "3.70"
Set-TextCell "E42" "  -8.51%  "
Set-TextCell "D43" "0.0327"
Set-TextCell "E43" "  -8.34%  "
Set-TextCell "D44" "13.61"
Set-TextCell "E44" "  +4.69%  "
Set-TextCell "E45" "  -0.07%  "
Set-TextCell "D46" "1.762.28"
Set-TextCell "E46" "  -4.99%  "
Set-TextCell "D47" "87.10"
Set-TextCell "E47" "  -10.92%  "
Set-TextCell "D48" "0.207"
Set-TextCell "E48" "  -10.00%  "
Set-TextCell "D49" "5.36"
Set-TextCell "E49" "  -10.30%  "
Set-TextCell "D50" "76.35"
Set-TextCell "E50" "  -8.74%  "
Set-TextCell "D51" "59.47"
Set-TextCell "E51" "  -15.68%  "
